$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.706.30'
$ws.Range("E2").Value = '  +1.48%  '
$ws.Range("D3").Value = '3.032.06'
$ws.Range("E3").Value = '  +2.84%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '380.83'
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("D6").Value = '103.16'
$ws.Range("E6").Value = '  +1.61%  '
$ws.Range("D7").Value = '0.546'
$ws.Range("E7").Value = '  +0.93%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.595'
$ws.Range("E9").Value = '  +2.27%  '
$ws.Range("D10").Value = '36.91'
$ws.Range("E10").Value = '  +2.02%  '
$ws.Range("E11").Value = '  -0.21%  '
$ws.Range("D12").Value = '0.0861'
$ws.Range("E12").Value = '  +1.49%  '
$ws.Range("D13").Value = '3.528.03'
$ws.Range("E13").Value = '  +3.37%  '
$ws.Range("D14").Value = '18.59'
$ws.Range("E14").Value = '  +1.66%  '
$ws.Range("D15").Value = '7.75'
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("D16").Value = '3.030.03'
$ws.Range("E16").Value = '  +2.67%  '
$ws.Range("E17").Value = '  -2.81%  '
$ws.Range("D18").Value = '10.48'
$ws.Range("E18").Value = '  -13.35%  '
$ws.Range("D19").Value = '51.727.40'
$ws.Range("E19").Value = '  +1.60%  '
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").Value = '12.50'
$ws.Range("E21").Value = '  +0.98%  '
$ws.Range("D22").Value = '0.0₃0962'
$ws.Range("E22").Value = '  +1.14%  '
$ws.Range("E23").Value = '  +1.14%  '
$ws.Range("D24").Value = '268.68'
$ws.Range("E24").Value = '  +0.96%  '
$ws.Range("D25").Value = '3.18'
$ws.Range("E25").Value = '  -0.23%  '
$ws.Range("D26").Value = '8.27'
$ws.Range("E26").Value = '  +2.07%  '
$ws.Range("D27").Value = '7.57'
$ws.Range("E27").Value = '  +7.72%  '
$ws.Range("D28").Value = '0.174'
$ws.Range("E28").Value = '  +6.42%  '
$ws.Range("D29").Value = '26.29'
$ws.Range("E29").Value = '  +2.73%  '
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").Value = '0.108'
$ws.Range("E31").Value = '  +0.46%  '
$ws.Range("D32").Value = '10.29'
$ws.Range("E32").Value = '  +1.13%  '
$ws.Range("D33").Value = '34.19'
$ws.Range("E33").Value = '  +1.85%  '
$ws.Range("D34").Value = '50.49'
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("D36").Value = '0.0448'
$ws.Range("E36").Value = '  +3.84%  '
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").Value = '3.34'
$ws.Range("E38").Value = '  +7.62%  '
$ws.Range("E39").Value = '  +13.85%  '
$ws.Range("D40").Value = '17.08'
$ws.Range("E40").Value = '  +3.09%  '
$ws.Range("E41").Value = '  +2.89%  '
$ws.Range("D42").Value = '2.57'
$ws.Range("E42").Value = '  +2.20%  '
$ws.Range("D43").Value = '127.49'
$ws.Range("E43").Value = '  +6.66%  '
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("D45").Value = '3.78'
$ws.Range("E45").Value = '  +7.12%  '
$ws.Range("D46").Value = '21.84'
$ws.Range("E46").Value = '  +2.10%  '
$ws.Range("D47").Value = '2.10'
$ws.Range("E47").Value = '  +4.27%  '
$ws.Range("D48").Value = '2.39'
$ws.Range("E48").Value = '  +3.00%  '
$ws.Range("D49").Value = '2.036.18'
$ws.Range("E49").Value = '  +1.77%  '
$ws.Range("D50").Value = '3.334.69'
$ws.Range("E50").Value = '  +2.93%  '
$ws.Range("D51").Value = '0.0320'
$ws.Range("E51").Value = '  +0.92%  '
